$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.642.63"
$ws.Range("E2").Value = "  +2.12%  "

$ws.Range("D3").Value = "1.889.79"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "244.58"
$ws.Range("E5").Value = "  +1.14%  "

$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").Value = "0.4961"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "0.2956"
$ws.Range("E8").Value = "  +1.48%  "

$ws.Range("D9").Value = "0.06802"
$ws.Range("E9").Value = "  +3.38%  "

$ws.Range("D10").Value = "1.889.97"
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("D11").Value = "17.03"
$ws.Range("E11").Value = "  +2.01%  "

$ws.Range("D12").Value = "0.07303"
$ws.Range("E12").Value = "  +1.92%  "

$ws.Range("D13").Value = "90.79"
$ws.Range("E13").Value = "  +5.49%  "

$ws.Range("D14").Value = "5.060"
$ws.Range("E14").Value = "  +4.74%  "

$ws.Range("D15").Value = "0.6712"
$ws.Range("E15").Value = "  +1.51%  "

$ws.Range("D16").Value = "30.639.42"
$ws.Range("E16").Value = "  +2.17%  "

$ws.Range("D17").Value = "0.000007916"
$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("E18").Value = "  +0.25%  "

$ws.Range("D19").Value = "13.19"
$ws.Range("E19").Value = "  +3.93%  "

$ws.Range("D20").Value = "2.137.78"
$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").Value = "4.849"
$ws.Range("E22").Value = "  +2.15%  "

$ws.Range("D23").Value = "175.46"
$ws.Range("E23").Value = "  +30.92%  "

$ws.Range("D24").Value = "6.049"
$ws.Range("E24").Value = "  +8.25%  "

$ws.Range("D25").Value = "9.261"
$ws.Range("E25").Value = "  +1.93%  "

$ws.Range("D26").Value = "155.57"
$ws.Range("E26").Value = "  +3.45%  "

$ws.Range("D27").Value = "18.70"
$ws.Range("E27").Value = "  +12.00%  "

$ws.Range("D28").Value = "1.923"
$ws.Range("E28").Value = "  +1.04%  "

$ws.Range("D29").Value = "1.392"
$ws.Range("E29").Value = "  +1.43%  "

$ws.Range("D30").Value = "4.337"
$ws.Range("E30").Value = "  +4.25%  "

$ws.Range("D31").Value = "0.08902"
$ws.Range("E31").Value = "  +2.73%  "

$ws.Range("D32").Value = "4.020"
$ws.Range("E32").Value = "  +2.22%  "

$ws.Range("D33").Value = "0.05230"
$ws.Range("E33").Value = "  +4.57%  "

$ws.Range("D34").Value = "0.7385"
$ws.Range("E34").Value = "  +4.82%  "

$ws.Range("D35").Value = "1.132"
$ws.Range("E35").Value = "  +3.15%  "

$ws.Range("D36").Value = "2.673"
$ws.Range("E36").Value = "  +0.73%  "

$ws.Range("D37").Value = "0.01874"
$ws.Range("E37").Value = "  +10.72%  "

$ws.Range("D38").Value = "2.706"
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("D39").Value = "2.172"
$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("D40").Value = "0.9341"
$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("D41").Value = "0.4357"
$ws.Range("E41").Value = "  +4.10%  "

$ws.Range("D42").Value = "105.70"
$ws.Range("E42").Value = "  +4.03%  "

$ws.Range("D43").Value = "5.805"
$ws.Range("E43").Value = "  -2.52%  "

$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("D45").Value = "7.646"
$ws.Range("E45").Value = "  +3.41%  "

$ws.Range("E46").Value = "  +8.00%  "

$ws.Range("D47").Value = "0.05825"
$ws.Range("E47").Value = "  +2.78%  "

$ws.Range("D48").Value = "33.41"
$ws.Range("E48").Value = "  +3.14%  "

$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "0.3872"
$ws.Range("E49").Value = "  +5.07%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "8.473"
$ws.Range("E50").Value = "  +5.37%  "

$ws.Range("D51").Value = "1.378"
$ws.Range("E51").Value = "  +3.52%  "
